# Apply the updated odds values (FlashScore 2024-10-21 refresh) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 odds updates
$ws.Range("M2").Value = 1.08
$ws.Range("O2").Value = 1.4
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8

# Row 3 odds updates
$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.4
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 1.8
$ws.Range("X3").Value = 11
$ws.Range("AJ3").Value = 12

# Row 4 odds updates
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.75
$ws.Range("AC4").Value = 8
$ws.Range("AG4").Value = 451
$ws.Range("AO4").Value = 11
$ws.Range("AZ4").Value = 81

# Row 5 odds updates
$ws.Range("G5").Value = 2.55
$ws.Range("H5").Value = 2.9
$ws.Range("Z5").Value = 23
$ws.Range("AA5").Value = 23
$ws.Range("AC5").Value = 6.5
$ws.Range("AD5").Value = 6
$ws.Range("AS5").Value = 251

# Row 23 odds updates
$ws.Range("G23").Value = 2.55
$ws.Range("I23").Value = 2.8
$ws.Range("J23").Value = 3.4
$ws.Range("K23").Value = 1.95
$ws.Range("L23").Value = 3.6
$ws.Range("M23").Value = 1.08
$ws.Range("O23").Value = 1.44
$ws.Range("P23").Value = 2.63
$ws.Range("S23").Value = 1.53
$ws.Range("T23").Value = 2.38
$ws.Range("V23").Value = 1.73
$ws.Range("W23").Value = 7
$ws.Range("X23").Value = 11
$ws.Range("AF23").Value = 67
$ws.Range("AI23").Value = 13
$ws.Range("AK23").Value = 29
$ws.Range("AL23").Value = 26
$ws.Range("AO23").Value = 15
$ws.Range("AT23").Value = 2.38
$ws.Range("AW23").Value = 4.75

# Row 24 odds updates
$ws.Range("G24").Value = 2.88
$ws.Range("I24").Value = 2.45
$ws.Range("L24").Value = 3.25
$ws.Range("U24").Value = 2
$ws.Range("V24").Value = 1.73
$ws.Range("AK24").Value = 23
$ws.Range("AW24").Value = 4.33
$ws.Range("AY24").Value = 26

# Row 39 odds updates
$ws.Range("G39").Value = 2.15
$ws.Range("I39").Value = 3.2
$ws.Range("Q39").Value = 1.98
$ws.Range("R39").Value = 1.88
$ws.Range("AC39").Value = 10
$ws.Range("AE39").Value = 15
$ws.Range("AJ39").Value = 12
$ws.Range("AL39").Value = 26
$ws.Range("AU39").Value = 8
$ws.Range("BA39").Value = 81
